$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns BC and BD on row 1
$ws.Range("BC1").Value = "Odd_CS_3-3_HT"
$ws.Range("BD1").Value = "Odd_CS_4-4_HT"

# Copy style from an existing header cell (e.g. BB1) to the new header cells
$ws.Range("BB1").Copy()
$ws.Range("BC1:BD1").PasteSpecial(-4122) # xlPasteFormats

# Update row 2 values
$ws.Range("A2").Value = "Ukw2Lkbe"
$ws.Range("C2").Value = "14:30"
$ws.Range("D2").Value = "AUSTRIA - 2. LIGA"
$ws.Range("E2").Value = "A. Lustenau"
$ws.Range("F2").Value = "Stripfing"
$ws.Range("G2").Value = 1.65
$ws.Range("H2").Value = 3.75
$ws.Range("I2").Value = 4.75
$ws.Range("J2").Value = 2.27
$ws.Range("K2").Value = 2.15
$ws.Range("L2").Value = 4.9
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 7.4
$ws.Range("O2").Value = 1.3
$ws.Range("P2").Value = 3.2
$ws.Range("Q2").Value = 1.88
$ws.Range("R2").Value = 1.82
$ws.Range("S2").Value = 1.42
$ws.Range("T2").Value = 2.67
$ws.Range("U2").Value = 1.88
$ws.Range("V2").Value = 1.83
$ws.Range("W2").Value = 6.6
$ws.Range("X2").Value = 7.5
$ws.Range("Y2").Value = 8.25
$ws.Range("Z2").Value = 12.5
$ws.Range("AA2").Value = 13.5
$ws.Range("AB2").Value = 28
$ws.Range("AC2").Value = 7.4
$ws.Range("AD2").Value = 7.2
$ws.Range("AE2").Value = 16.5
$ws.Range("AF2").Value = 80
$ws.Range("AG2").Value = 700
$ws.Range("AH2").Value = 12.5
$ws.Range("AI2").Value = 26
$ws.Range("AJ2").Value = 15.5
$ws.Range("AK2").Value = 80
$ws.Range("AL2").Value = 45
$ws.Range("AM2").Value = 50
$ws.Range("AN2").Value = 3.45
$ws.Range("AO2").Value = 8.25
$ws.Range("AP2").Value = 19
$ws.Range("AQ2").Value = 28
$ws.Range("AR2").Value = 65
$ws.Range("AS2").Value = 300
$ws.Range("AT2").Value = 2.67
$ws.Range("AU2").Value = 7.9
$ws.Range("AV2").Value = 80
$ws.Range("AW2").Value = 6.3
$ws.Range("AX2").Value = 27
$ws.Range("AY2").Value = 35
$ws.Range("AZ2").Value = 175
$ws.Range("BA2").Value = 200
$ws.Range("BB2").Value = 500
$ws.Range("BC2").Value = 81
$ws.Range("BD2").Value = 81
